$d = $word.ActiveDocument

# --- Step 1: merge runs at the end of paragraph 78 ("26.   " + "Object linkage/Prototypal inheritance") ---
$p78 = $d.Paragraphs.Item(78)
$p78Start = $p78.Range.Start
$p78End = $p78.Range.End
$mergeRng1 = $d.Range($p78Start + 1, $p78End)
$mergeRng1.Text = "26.   Object linkage/Prototypal inheritance"

# --- Step 2: merge runs at the end of paragraph 80 ("Object linkage" + "… and Functions") ---
$p80 = $d.Paragraphs.Item(80)
$p80Start = $p80.Range.Start
$p80End = $p80.Range.End
$mergeRng2 = $d.Range($p80Start + 45, $p80End)
$mergeRng2.Text = "Object linkage… and Functions"

# --- Step 3: apply strikethrough to the affected paragraphs (struck-through / addressed items) ---
$strikeParagraphs = @(6,7,8,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,93,94,95)
foreach ($idx in $strikeParagraphs) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.StrikeThrough = 1
}

# --- Step 4: apply strikethrough to the hyperlink runs inside those paragraphs ---
# (hyperlink field runs are not reached by Paragraph.Range.Font, so set them directly)
$strikeHyperlinks = @(1,2,25,26,27,28,29,30,31,32)
foreach ($idx in $strikeHyperlinks) {
    $h = $d.Hyperlinks.Item($idx)
    $h.Range.Font.StrikeThrough = 1
}

Write-Output "strikethrough applied"
